# Weekly update: insert two new daily-price rows (row 378 and 379) above the
# existing data, shifting the rest of the table down by two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 378:379 (shifts 378:486 down to 380:488,
# carrying the existing cell formatting, e.g. the date style on column D).
$ws.Range("A378:R379").Insert()

# --- New row 378 ---
$ws.Cells.Item(378, 1).Value = 6
$ws.Cells.Item(378, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(378, 3).Value = "Metropolitana"
$ws.Cells.Item(378, 4).Value = 44463
$ws.Cells.Item(378, 5).Value = 13
$ws.Cells.Item(378, 6).Value = 100112009
$ws.Cells.Item(378, 7).Value = "Acelga"
$ws.Cells.Item(378, 8).Value = "Sin especificar"
$ws.Cells.Item(378, 9).Value = "Primera"
$ws.Cells.Item(378, 10).Value = 280
$ws.Cells.Item(378, 11).Value = 8000
$ws.Cells.Item(378, 12).Value = 10000
$ws.Cells.Item(378, 13).Value = 9143
$ws.Cells.Item(378, 14).Value = "`$/docena de atados"
$ws.Cells.Item(378, 15).Value = "Región Metropolitana"
$ws.Cells.Item(378, 16).Value = 3048
$ws.Cells.Item(378, 17).Value = 3
$ws.Cells.Item(378, 18).Value = "Hortaliza"

# --- New row 379 ---
$ws.Cells.Item(379, 1).Value = 6
$ws.Cells.Item(379, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(379, 3).Value = "Metropolitana"
$ws.Cells.Item(379, 4).Value = 44463
$ws.Cells.Item(379, 5).Value = 13
$ws.Cells.Item(379, 6).Value = 100112009
$ws.Cells.Item(379, 7).Value = "Acelga"
$ws.Cells.Item(379, 8).Value = "Sin especificar"
$ws.Cells.Item(379, 9).Value = "Segunda"
$ws.Cells.Item(379, 10).Value = 40
$ws.Cells.Item(379, 11).Value = 7000
$ws.Cells.Item(379, 12).Value = 7000
$ws.Cells.Item(379, 13).Value = 7000
$ws.Cells.Item(379, 14).Value = "`$/docena de atados"
$ws.Cells.Item(379, 15).Value = "Región Metropolitana"
$ws.Cells.Item(379, 16).Value = 2333
$ws.Cells.Item(379, 17).Value = 3
$ws.Cells.Item(379, 18).Value = "Hortaliza"
